# Apply the "Monday Aug 25 22:05" edit:
#  - Collapse the day table from 6 data rows (Day 4..Day 9) down to a single
#    data row: A2 = 1 (number), B2 = new concat path.
#  - Delete the now-unused rows 3-7.
#  - Move the active selection to A3 (first empty row beneath the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 7 (the old Day 5..Day 9 rows), leaving just the
# header (row 1) and the single data row (row 2).
$ws.Range("A3:B7").EntireRow.Delete() | Out-Null

# Row 2: A2 becomes a plain number, B2 becomes the new path string.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "/Users/miguel_barron/Downloads/concatEphys/03312025"

# Match the saved selection state from the diff.
$ws.Range("A3").Select() | Out-Null
